$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row re-orderings)
# Each value is prefixed with a literal apostrophe so Excel stores it as
# plain text (matching the source inlineStr cells) instead of coercing
# numeric-looking strings (e.g. "1.00", "0.492") into floating point numbers.

$ws.Range("D2").Value = "'64.852.82"
$ws.Range("E2").Value = "'  -0.40%  "
$ws.Range("D3").Value = "'3.557.54"
$ws.Range("E3").Value = "'  +2.34%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'597.97"
$ws.Range("E5").Value = "'  +1.76%  "
$ws.Range("D6").Value = "'134.83"
$ws.Range("E6").Value = "'  -1.49%  "
$ws.Range("D7").Value = "'3.555.75"
$ws.Range("E7").Value = "'  +2.30%  "
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "'  +0.33%  "
$ws.Range("E10").Value = "'  +0.10%  "
$ws.Range("D11").Value = "'6.94"
$ws.Range("E11").Value = "'  -3.12%  "
$ws.Range("D12").Value = "'0.383"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("D13").Value = "'4.160.94"
$ws.Range("E13").Value = "'  +2.30%  "
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'3.561.44"
$ws.Range("E14").Value = "'  +2.19%  "
$ws.Range("B15").Value = "'ShibaInu"
$ws.Range("C15").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000180"
$ws.Range("E15").Value = "'  -0.43%  "
$ws.Range("D16").Value = "'26.82"
$ws.Range("E16").Value = "'  +0.97%  "
$ws.Range("E17").Value = "'  +0.42%  "
$ws.Range("D18").Value = "'65.002.04"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("D19").Value = "'9.94"
$ws.Range("E19").Value = "'  +2.48%  "
$ws.Range("D20").Value = "'14.27"
$ws.Range("E20").Value = "'  +2.48%  "
$ws.Range("D21").Value = "'5.79"
$ws.Range("E21").Value = "'  +0.34%  "
$ws.Range("D22").Value = "'387.31"
$ws.Range("E22").Value = "'  -0.26%  "
$ws.Range("D23").Value = "'0.575"
$ws.Range("E23").Value = "'  +3.51%  "
$ws.Range("D24").Value = "'3.703.75"
$ws.Range("E24").Value = "'  +2.36%  "
$ws.Range("D25").Value = "'73.76"
$ws.Range("E25").Value = "'  +1.78%  "
$ws.Range("E26").Value = "'  +0.04%  "
$ws.Range("D27").Value = "'0.0000113"
$ws.Range("E27").Value = "'  +2.80%  "
$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "'  +3.98%  "
$ws.Range("E29").Value = "'  +0.17%  "
$ws.Range("E30").Value = "'  +2.99%  "
$ws.Range("D31").Value = "'8.35"
$ws.Range("E31").Value = "'  +2.78%  "
$ws.Range("D32").Value = "'1.47"
$ws.Range("E32").Value = "'  +23.39%  "
$ws.Range("D33").Value = "'3.559.14"
$ws.Range("E33").Value = "'  +1.78%  "
$ws.Range("B34").Value = "'USDe"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.02%  "
$ws.Range("B35").Value = "'EthereumClassic"
$ws.Range("C35").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'23.88"
$ws.Range("E35").Value = "'  +3.70%  "
$ws.Range("D36").Value = "'0.143"
$ws.Range("E36").Value = "'  +0.04%  "
$ws.Range("B37").Value = "'Monero"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'168.29"
$ws.Range("E37").Value = "'  -1.27%  "
$ws.Range("B38").Value = "'Aptos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'6.87"
$ws.Range("E38").Value = "'  +0.85%  "
$ws.Range("D39").Value = "'1.53"
$ws.Range("E39").Value = "'  +3.67%  "
$ws.Range("D40").Value = "'4.94"
$ws.Range("E40").Value = "'  +4.75%  "
$ws.Range("D41").Value = "'0.0799"
$ws.Range("E41").Value = "'  +2.90%  "
$ws.Range("D42").Value = "'0.822"
$ws.Range("E42").Value = "'  +1.53%  "
$ws.Range("D43").Value = "'26.67"
$ws.Range("E43").Value = "'  +6.65%  "
$ws.Range("D44").Value = "'42.63"
$ws.Range("E44").Value = "'  +0.37%  "
$ws.Range("E45").Value = "'  +0.00%  "
$ws.Range("D46").Value = "'4.42"
$ws.Range("E46").Value = "'  +1.89%  "
$ws.Range("E47").Value = "'  +3.76%  "
$ws.Range("D48").Value = "'1.63"
$ws.Range("E48").Value = "'  +0.57%  "
$ws.Range("D49").Value = "'2.475.36"
$ws.Range("E49").Value = "'  +11.74%  "
$ws.Range("D50").Value = "'6.86"
$ws.Range("E50").Value = "'  +2.62%  "
$ws.Range("D51").Value = "'0.861"
$ws.Range("E51").Value = "'  +7.69%  "
